# Applies the "Added tests for verifying calculations on real map" commit:
#  - Inserts two extra perimeter points into the existing POS1..POS18 grid
#    (becoming POS1..POS21), which re-numbers/shifts Q/R (and therefore S/T/U)
#    for rows 11-24 and adds new rows 25-27.
#  - Adds a block of "real map" verification rows (31-54) that compare the
#    computed lat/long strings against real aircraft lat/long fixes, plus
#    some blank formatted rows (55-60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the Q/R offsets for the existing perimeter rows (7-24). P and the
#    formulas in S/T/U stay the same text - they just recompute from Q/R.
#    Rows 7-10 are unchanged.
# ---------------------------------------------------------------------------

$ws.Range("R11").Value = 0
$ws.Range("R12").Value = 2
$ws.Range("Q13").Value = -6
$ws.Range("Q14").Value = -4
$ws.Range("Q15").Value = -2
$ws.Range("Q16").Value = 0
$ws.Range("Q17").Value = 2
$ws.Range("Q18").Value = 4
$ws.Range("R18").Value = 4
$ws.Range("R19").Value = 4
$ws.Range("R20").Value = 2
$ws.Range("Q21").Value = 6
$ws.Range("R21").Value = 0
$ws.Range("Q22").Value = 6
$ws.Range("R22").Value = -2
$ws.Range("Q23").Value = 6
$ws.Range("R23").Value = -4
$ws.Range("Q24").Value = 4
$ws.Range("R24").Value = -4

# Cells S11,T11,U11 / S16,T16,U16 / S21,T21,U21 / S26? (no) had their shared
# formula group broken in the real edit (Excel re-derives new si's whenever a
# formula is re-entered). Re-enter them explicitly so the cached values match;
# also restore "Normal" style afterwards since re-entering a formula on this
# engine otherwise guesses a stray number format from a neighbouring cell.
function Set-Formula($ref, $formula) {
    $ws.Range($ref).Formula = $formula
    $ws.Range($ref).Style = "Normal"
}

Set-Formula "S11" '=S$6+Q11*S$4+0.01'
Set-Formula "T11" '=T$6+R11*T$4+0.01'
$ws.Range("U11").Formula = '=CONCAT("((",S11,",",T11,"), ""POS",P11,"""),\")'

Set-Formula "S16" '=S$6+Q16*S$4+0.01'
Set-Formula "T16" '=T$6+R16*T$4+0.01'
$ws.Range("U16").Formula = '=CONCAT("((",S16,",",T16,"), ""POS",P16,"""),\")'

Set-Formula "S21" '=S$6+Q21*S$4+0.01'
Set-Formula "T21" '=T$6+R21*T$4+0.01'
$ws.Range("U21").Formula = '=CONCAT("((",S21,",",T21,"), ""POS",P21,"""),\")'

# ---------------------------------------------------------------------------
# 2. Add the 3 new perimeter rows 25-27 (POS19, POS20, POS21) continuing the
#    pattern established by rows 7-24.
# ---------------------------------------------------------------------------

$newRows = @(
    @{ Row = 25; P = 19; Q = 2;  R = -4 },
    @{ Row = 26; P = 20; Q = 2;  R = -6 },
    @{ Row = 27; P = 21; Q = 2;  R = -8 }
)
foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Formula = "=S`$6+Q$row*S`$4+0.01"
    $ws.Range("S$row").Style = "Normal"
    $ws.Range("T$row").Formula = "=T`$6+R$row*T`$4+0.01"
    $ws.Range("T$row").Style = "Normal"
    $ws.Range("U$row").Formula = "=CONCAT(`"((`",S$row,`",`",T$row,`"), `"`"POS`",P$row,`"`"`"),\`")"
}

# re-apply the P/Q/R border+center style used by the rest of the grid
$ws.Range("P25:R27").Style = $ws.Range("P24").Style

# ---------------------------------------------------------------------------
# 3. "Real map" verification block (rows 31-54).
# ---------------------------------------------------------------------------

$ws.Range("S31").Formula = '=CONCAT(T6,",",S6)'

for ($i = 7; $i -le 27; $i++) {
    $targetRow = $i + 25
    $ws.Range("S$targetRow").Formula = "=CONCAT(T$i,`",`",S$i)"
}

# Tail / nose numbers and timestamps of the real aircraft used to validate
# the grid (columns C/F/G/H), set in the same order the author typed them so
# the shared-string table is built up with matching indices.
$ws.Range("C33").NumberFormat = "[h]:mm:ss"
$ws.Range("C33").Value = 4.291666666666667
$ws.Range("F33").Value = 33.2746
$ws.Range("G33").Value = -111.8021
$ws.Range("H33").Value = "N7098P"

$ws.Range("C34").NumberFormat = "[h]:mm:ss"
$ws.Range("C34").Value = 4.333333333333333
$ws.Range("F34").Value = 33.251300000000001
$ws.Range("G34").Value = -111.80459999999999
$ws.Range("H34").Value = "MSQT818"

$ws.Range("C35").NumberFormat = "[h]:mm:ss"
$ws.Range("C35").Value = 25.041666666666668
$ws.Range("F35").Value = 33.335099999999997
$ws.Range("G35").Value = -111.6497
$ws.Range("H35").Value = "N14053"

$ws.Range("C36").NumberFormat = "[h]:mm:ss"
$ws.Range("C36").Value = 25.125
$ws.Range("F36").Value = 33.288400000000003
$ws.Range("G36").Value = -111.6323
$ws.Range("H36").Value = "MSQT955"

$ws.Range("C37").NumberFormat = "[h]:mm:ss"
$ws.Range("C37").Value = 29.291666666666668
$ws.Range("F37").Value = 33.2926
$ws.Range("G37").Value = -111.6193
$ws.Range("H37").Value = "NDU531"

$ws.Range("C38").NumberFormat = "[h]:mm:ss"
$ws.Range("C38").Value = 29.333333333333332
$ws.Range("F38").Value = 33.256700000000002
$ws.Range("G38").Value = -111.5945
$ws.Range("H38").Value = "N94HL"

# Rows 47-54: a little lat/long table ("(", ",", ")" used so the
# CONCAT(...) results read like "(lat, lon)" literals), sorted by tail
# number (hence the trailing <sortState> on D47:F58 in the original file).
$rows4754 = @(
    @{ Row = 47; Tail = "*******"; F = 33.331400000000002;  H = -111.7671 },
    @{ Row = 48; Tail = "*******"; F = 33.319299999999998;  H = -111.7709 },
    @{ Row = 49; Tail = "MSQT182"; F = 33.2821;             H = -111.79130000000001 },
    @{ Row = 50; Tail = "MSQT182"; F = 33.279400000000003;  H = -111.795 },
    @{ Row = 51; Tail = "N4400Q";  F = 33.281199999999998;  H = -111.623 },
    @{ Row = 52; Tail = "N4400Q";  F = 33.319299999999998;  H = -111.6677 },
    @{ Row = 53; Tail = "N7098P";  F = 33.277999999999999;  H = -111.79559999999999 },
    @{ Row = 54; Tail = "N7098P";  F = 33.273800000000001;  H = -111.8018 }
)

foreach ($r in $rows4754) {
    $row = $r.Row
    $ws.Range("C$row").NumberFormat = "[h]:mm:ss"
    $ws.Range("D$row").Value = $r.Tail
    $ws.Range("E$row").Value = "("
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = ","
    $ws.Range("H$row").Value = $r.H
    $ws.Range("J$row").Value = ")"
    $ws.Range("K$row").Formula = "=CONCAT(F${row}:H$row)"
}

# S53 breaks from the S32:S52 shared-formula pattern, comparing back to the
# first real fix (row 33) instead of to the synthetic grid.
$ws.Range("S53").Formula = '=CONCAT(F33,",",G33)'

# Trailing blank rows (55-60), still carrying the time number format so the
# column reads consistently all the way down.
for ($row = 55; $row -le 60; $row++) {
    $ws.Range("C$row").NumberFormat = "[h]:mm:ss"
}

# ---------------------------------------------------------------------------
# 4. View state: scroll/selection moved down to the new data while editing.
# ---------------------------------------------------------------------------
$ws.Range("N60").Select()
